{"js": "// 1. Update the date paragraph from 2020-08-03 to 2020-10-08.\nconst dateResults = context.document.body.search(\"2020-08-03\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"2020-10-08\", \"Replace\");\n}\nawait context.sync();\n\n// 2. Append a new \"Compact\" list item (same numbered list as the surrounding\n//    links, numId 1003) with a hyperlink to \"Happy Git and GitHub for the useR\",\n//    right after the existing \"Bayesian inference with INLA\" entry at the end\n//    of the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Happy Git and GitHub for the useR\",\n  \"After\"\n);\n\nconst newRange = newParagraph.getRange();\nnewRange.hyperlink = \"https://happygitwithr.com/\";\n\nawait context.sync();\n", "ps1": "# 1. Update the date paragraph from 2020-08-03 to 2020-10-08.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Date\") {\n        $p.Range.Text = \"2020-10-08\"\n        break\n    }\n}\n\n# 2. Append a new \"Compact\" list item (same numbered list as the surrounding\n#    links, numId 1003) with a hyperlink to \"Happy Git and GitHub for the useR\",\n#    right after the existing \"Bayesian inference with INLA\" entry at the end\n#    of the document body.\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newRange = $newParagraph.Range\n$newRange.InsertAfter(\"Happy Git and GitHub for the useR\")\n\n$newParagraph2 = $d.Paragraphs.Last\n$newParagraph2.Range.Hyperlink = \"https://happygitwithr.com/\"\n"}
